# VISA Marriott 2017.xlsx -- append newly-imported transaction rows to Sheet1.
#
# The sheet is a running list of card transactions grouped into "statement"
# blocks separated by blank rows. This edit appends one more statement block
# (rows 34-43) below the existing data (which ended at row 33), following the
# exact same layout: Type | Trans Date | Post Date | Description | Amount,
# with an occasional helper formula in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlPasteFormats - used to replicate the look (number format / font / fill)
# of an existing data row onto the freshly-appended ones, the same way the
# original author likely did (typing values into a copy of the row above).
$xlPasteFormats = -4122

function Copy-RowFormat {
    param($Row, $TemplateRow)
    $ws.Range("A$TemplateRow`:E$TemplateRow").Copy() | Out-Null
    $ws.Range("A$Row`:E$Row").PasteSpecial($xlPasteFormats) | Out-Null
}

# Stamp the formatting (number formats / fonts) for every new row first --
# PasteSpecial-formats never touches the shared-string table, so this can
# happen in any order.
Copy-RowFormat 34 33
Copy-RowFormat 35 33
Copy-RowFormat 37 33
Copy-RowFormat 38 33
Copy-RowFormat 39 33
Copy-RowFormat 40 33
Copy-RowFormat 41 33
Copy-RowFormat 42 33
Copy-RowFormat 43 33

$ws.Range("B31:C31").Copy() | Out-Null
$ws.Range("B36:C36").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("E40").Copy() | Out-Null
$ws.Range("F40").PasteSpecial($xlPasteFormats) | Out-Null

# Fill in the merchant descriptions bottom-up, i.e. in the order the rows
# were actually transcribed off the statement, newest first.
$ws.Cells.Item(43, 4).Value = "EXOTIC THAI CAFE"
$ws.Cells.Item(42, 4).Value = "EAST COAST PIZZA CO"
$ws.Cells.Item(41, 4).Value = "VAN NUYS AIRPORT PARKING"
$ws.Cells.Item(40, 4).Value = "STARBUCKS MARINA CABOS"
$ws.Cells.Item(39, 4).Value = "FUR"
$ws.Cells.Item(37, 4).Value = "WWW COSTCO COM"
$ws.Cells.Item(34, 4).Value = "PRESTIGE MEDICAL ASSOCIA"
$ws.Cells.Item(35, 4).Value = "VITALITY APPLEWATCH"
$ws.Cells.Item(38, 4).Value = "AUTOMATIC PAYMENT - THANK"

# --- row 34 ------------------------------------------------------------
$ws.Cells.Item(34, 1).Value = "Sale"
$ws.Cells.Item(34, 2).Value = 42838
$ws.Cells.Item(34, 3).Value = 42839
$ws.Cells.Item(34, 5).Value = -31.84

# --- row 35 ------------------------------------------------------------
$ws.Cells.Item(35, 1).Value = "Sale"
$ws.Cells.Item(35, 2).Value = 42842
$ws.Cells.Item(35, 3).Value = 42843
$ws.Cells.Item(35, 5).Value = -4

# --- row 37 ------------------------------------------------------------
$ws.Cells.Item(37, 1).Value = "Sale"
$ws.Cells.Item(37, 2).Value = 42857
$ws.Cells.Item(37, 3).Value = 42858
$ws.Cells.Item(37, 5).Value = -70.74

# --- row 38 ------------------------------------------------------------
$ws.Cells.Item(38, 1).Value = "Payment"
$ws.Cells.Item(38, 2).Value = 42860
$ws.Cells.Item(38, 3).Value = 42860
$ws.Cells.Item(38, 5).Value = 156.46

# --- row 39 ------------------------------------------------------------
$ws.Cells.Item(39, 1).Value = "Sale"
$ws.Cells.Item(39, 2).Value = 42860
$ws.Cells.Item(39, 3).Value = 42863
$ws.Cells.Item(39, 5).Value = -77

# --- row 40 (has a helper formula in column F) --------------------------
$ws.Cells.Item(40, 1).Value = "Sale"
$ws.Cells.Item(40, 2).Value = 42866
$ws.Cells.Item(40, 3).Value = 42869
$ws.Cells.Item(40, 5).Value = -22.43
$ws.Cells.Item(40, 6).Formula = "=420/E40"

# --- row 41 ------------------------------------------------------------
$ws.Cells.Item(41, 1).Value = "Sale"
$ws.Cells.Item(41, 2).Value = 42869
$ws.Cells.Item(41, 3).Value = 42870
$ws.Cells.Item(41, 5).Value = -28

# --- row 42 ------------------------------------------------------------
$ws.Cells.Item(42, 1).Value = "Sale"
$ws.Cells.Item(42, 2).Value = 42869
$ws.Cells.Item(42, 3).Value = 42871
$ws.Cells.Item(42, 5).Value = -9.51

# --- row 43 ------------------------------------------------------------
$ws.Cells.Item(43, 1).Value = "Sale"
$ws.Cells.Item(43, 2).Value = 42870
$ws.Cells.Item(43, 3).Value = 42871
$ws.Cells.Item(43, 5).Value = -11.73

# Match the page's print orientation (added when the workbook was reprinted).
$ws.PageSetup.Orientation = 1

# Reflect the author's final on-screen selection.
$ws.Range("C35").Select()
